# Auto-generated edit script applying numeric corrections to the Leve profit
# calculation columns (H-N) across all 8 class sheets (ALC, ARM, BSM, CRP, CUL,
# GSM, LTW, WVR), per the scheduled runner data refresh.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 219.71428
$ws.Range("I2").Value = 222.83333
$ws.Range("J2").Value = 201
$ws.Range("K2").Value = 222.83333
$ws.Range("L2").Value = 201
$ws.Range("M2").Value = -109.83333
$ws.Range("N2").Value = -427
$ws.Range("H43").Value = 1657.1428
$ws.Range("J43").Value = 2150
$ws.Range("L43").Value = 2150
$ws.Range("N43").Value = -2288
$ws.Range("H74").Value = 4875
$ws.Range("I74").Value = 4875
$ws.Range("K74").Value = 4875
$ws.Range("M74").Value = -3939
$ws.Range("H75").Value = 86332
$ws.Range("I75").Value = 140000
$ws.Range("K75").Value = 140000
$ws.Range("M75").Value = -139064
$ws.Range("H77").Value = 4875
$ws.Range("I77").Value = 4875
$ws.Range("K77").Value = 24375
$ws.Range("M77").Value = -19695
$ws.Range("H78").Value = 86332
$ws.Range("I78").Value = 140000
$ws.Range("K78").Value = 420000
$ws.Range("M78").Value = -415320
$ws.Range("H100").Value = 3200.2
$ws.Range("I100").Value = 3337.25
$ws.Range("K100").Value = 3337.25
$ws.Range("M100").Value = -2796.25
$ws.Range("H101").Value = 365.125
$ws.Range("I101").Value = 321.16666
$ws.Range("J101").Value = 497
$ws.Range("K101").Value = 963.4999799999999
$ws.Range("L101").Value = 1491
$ws.Range("M101").Value = 658.5000200000001
$ws.Range("N101").Value = -4735
$ws.Range("H112").Value = 2108.3333
$ws.Range("I112").Value = 1200
$ws.Range("J112").Value = 2411.111
$ws.Range("K112").Value = 3600
$ws.Range("L112").Value = 7233.333
$ws.Range("M112").Value = -2492
$ws.Range("N112").Value = -9449.332999999999
$ws.Range("H138").Value = 3043.2173
$ws.Range("I138").Value = 1374.25
$ws.Range("J138").Value = 3933.3333
$ws.Range("K138").Value = 4122.75
$ws.Range("L138").Value = 11799.9999
$ws.Range("M138").Value = 1017.25
$ws.Range("N138").Value = -22079.9999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2758.0625
$ws.Range("I61").Value = 2344
$ws.Range("J61").Value = 3669
$ws.Range("K61").Value = 2344
$ws.Range("L61").Value = 3669
$ws.Range("M61").Value = -2132
$ws.Range("N61").Value = -4093
$ws.Range("H122").Value = 3030.6667
$ws.Range("I122").Value = 2472
$ws.Range("K122").Value = 7416
$ws.Range("M122").Value = -4966
$ws.Range("H136").Value = 2758.0625
$ws.Range("I136").Value = 2344
$ws.Range("J136").Value = 3669
$ws.Range("K136").Value = 7032
$ws.Range("L136").Value = 11007
$ws.Range("M136").Value = -4482
$ws.Range("N136").Value = -16107

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 20314
$ws.Range("J76").Value = 20314
$ws.Range("L76").Value = 20314
$ws.Range("N76").Value = -20944
$ws.Range("H79").Value = 20314
$ws.Range("J79").Value = 20314
$ws.Range("L79").Value = 20314
$ws.Range("N79").Value = -22498
$ws.Range("H86").Value = 3043.9375
$ws.Range("I86").Value = 3022.5557
$ws.Range("J86").Value = 3071.4285
$ws.Range("K86").Value = 3022.5557
$ws.Range("L86").Value = 3071.4285
$ws.Range("M86").Value = -1899.5557
$ws.Range("N86").Value = -5317.4285
$ws.Range("H89").Value = 3043.9375
$ws.Range("I89").Value = 3022.5557
$ws.Range("J89").Value = 3071.4285
$ws.Range("K89").Value = 15112.7785
$ws.Range("L89").Value = 15357.1425
$ws.Range("M89").Value = -9496.7785
$ws.Range("N89").Value = -26589.1425
$ws.Range("H107").Value = 3251.875
$ws.Range("I107").Value = 2380.25
$ws.Range("K107").Value = 2380.25
$ws.Range("M107").Value = -460.25
$ws.Range("H134").Value = 7944.8
$ws.Range("I134").Value = 7256
$ws.Range("J134").Value = 10700
$ws.Range("K134").Value = 21768
$ws.Range("L134").Value = 32100
$ws.Range("M134").Value = -19233
$ws.Range("N134").Value = -37170

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 1067
$ws.Range("I45").Value = 1067
$ws.Range("K45").Value = 1067
$ws.Range("M45").Value = -474
$ws.Range("H58").Value = 1547.2727
$ws.Range("J58").Value = 1517.091
$ws.Range("L58").Value = 1517.091
$ws.Range("N58").Value = -1923.091
$ws.Range("H62").Value = 4999.5
$ws.Range("I62").Value = 4999.5
$ws.Range("K62").Value = 4999.5
$ws.Range("M62").Value = -4375.5
$ws.Range("H65").Value = 4999.5
$ws.Range("I65").Value = 4999.5
$ws.Range("K65").Value = 24997.5
$ws.Range("M65").Value = -21877.5
$ws.Range("H134").Value = 3355.111
$ws.Range("I134").Value = 3039.5715
$ws.Range("J134").Value = 4459.5
$ws.Range("K134").Value = 9118.7145
$ws.Range("L134").Value = 13378.5
$ws.Range("M134").Value = -6583.7145
$ws.Range("N134").Value = -18448.5
$ws.Range("H136").Value = 1547.2727
$ws.Range("J136").Value = 1517.091
$ws.Range("L136").Value = 4551.272999999999
$ws.Range("N136").Value = -9651.272999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 2412
$ws.Range("I6").Value = 309.83334
$ws.Range("K6").Value = 929.5000200000001
$ws.Range("M6").Value = -816.5000200000001
$ws.Range("H132").Value = 1499
$ws.Range("I132").Value = 1499
$ws.Range("K132").Value = 13491
$ws.Range("M132").Value = -10961

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 503864.84
$ws.Range("I70").Value = 603439.8
$ws.Range("J70").Value = 5990
$ws.Range("K70").Value = 603439.8
$ws.Range("L70").Value = 5990
$ws.Range("M70").Value = -603169.8
$ws.Range("N70").Value = -6530
$ws.Range("H73").Value = 503864.84
$ws.Range("I73").Value = 603439.8
$ws.Range("J73").Value = 5990
$ws.Range("K73").Value = 603439.8
$ws.Range("L73").Value = 5990
$ws.Range("M73").Value = -602503.8
$ws.Range("N73").Value = -7862
$ws.Range("H107").Value = 328.9
$ws.Range("I107").Value = 328.9
$ws.Range("K107").Value = 328.9
$ws.Range("M107").Value = 1591.1

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1916.8334
$ws.Range("I100").Value = 1916.8334
$ws.Range("K100").Value = 1916.8334
$ws.Range("M100").Value = -1375.8334
$ws.Range("H136").Value = 2813.375
$ws.Range("I136").Value = 2700.6
$ws.Range("J136").Value = 4505
$ws.Range("K136").Value = 8101.799999999999
$ws.Range("L136").Value = 13515
$ws.Range("M136").Value = -5551.799999999999
$ws.Range("N136").Value = -18615

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 30000
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H81").Value = 3124.75
$ws.Range("I81").Value = 3000
$ws.Range("K81").Value = 6000
$ws.Range("M81").Value = -4939
$ws.Range("H84").Value = 3124.75
$ws.Range("I84").Value = 3000
$ws.Range("K84").Value = 30000
$ws.Range("M84").Value = -24696
$ws.Range("H122").Value = 4280.8
$ws.Range("J122").Value = 4751
$ws.Range("L122").Value = 14253
$ws.Range("N122").Value = -19153
$ws.Range("H132").Value = 1566.3334
$ws.Range("I132").Value = 1566.3334
$ws.Range("K132").Value = 4699.0002
$ws.Range("M132").Value = -2169.0002
$ws.Range("H136").Value = 3501.6924
$ws.Range("I136").Value = 3415
$ws.Range("J136").Value = 3665.4443
$ws.Range("K136").Value = 10245
$ws.Range("L136").Value = 10996.3329
$ws.Range("M136").Value = -7695
$ws.Range("N136").Value = -16096.3329
